$d = $word.ActiveDocument

# The results table's header row reads: Variable | Group1 | Group2 | Overall.
# Rename the third header cell from "Group2" to "Group0".
$table = $d.Tables(1)
$cell = $table.Cell(1, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null   # wdCharacter: drop the trailing end-of-cell mark

if ($cellRange.Text -eq "Group2") {
    $cellRange.Text = "Group0"
} else {
    # Fallback in case the table layout ever differs from what we expect.
    $d.Content.Find.Execute("Group2", $true, $true, $false, $false, $false, $true, 1, $false, "Group0", 2) | Out-Null
}
